$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column H corresponds to "Houses of Worship" policy indicator.
# Rows 33-176 had value 1 and are being updated to 0.
$ws.Range("H33:H176").Value = 0
